# Update ECO Actual (B) and ECO Balance (D) figures for each FSR row.
# Source values are stored as text (e.g. "9.00"), so a leading apostrophe
# is used to force Excel to keep the new values as text as well.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'7.00"
$ws.Range("D2").Value = "'7.00"

$ws.Range("B3").Value = "'15.00"
$ws.Range("D3").Value = "'15.00"

$ws.Range("B4").Value = "'7.00"
$ws.Range("D4").Value = "'7.00"

$ws.Range("B5").Value = "'11.00"
$ws.Range("D5").Value = "'11.00"

$ws.Range("B6").Value = "'4.00"
$ws.Range("D6").Value = "'4.00"

$ws.Range("B7").Value = "'44.00"
$ws.Range("D7").Value = "'44.00"
